# Orders workbook refresh: replace sample/test order rows with new data
# and drop the trailing (8th) order row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 8) entirely - it no longer exists after the edit.
$ws.Rows(8).Delete()

# Row 2
$ws.Cells.Item(2, 2).Value = "dgfsdgdf"
$ws.Cells.Item(2, 3).Value = "'+998908171355"
$ws.Cells.Item(2, 4).Value = "STORE_OWNER: sdfgdsfg"
$ws.Cells.Item(2, 5).Value = "DELIVERED"
$ws.Cells.Item(2, 7).Value = 12222
$ws.Cells.Item(2, 8).Value = 44881.44031289352
$ws.Cells.Item(2, 9).Value = 44881.443963356476

# Row 3
$ws.Cells.Item(3, 2).Value = "fdsgsdfg"
$ws.Cells.Item(3, 3).Value = "'+998908171355"
$ws.Cells.Item(3, 4).Value = "STORE_OWNER: fdsdfsdf"
$ws.Cells.Item(3, 5).Value = "DELIVERED"
$ws.Cells.Item(3, 7).Value = 10000
$ws.Cells.Item(3, 8).Value = 44881.449103900464
$ws.Cells.Item(3, 9).Value = 44881.450163761576

# Row 4
$ws.Cells.Item(4, 2).Value = "sdfgsdf"
$ws.Cells.Item(4, 3).Value = "'+998908171355"
$ws.Cells.Item(4, 4).Value = "STORE_OWNER: bjkbbbhjs"
$ws.Cells.Item(4, 5).Value = "DELIVERED"
$ws.Cells.Item(4, 7).Value = 122219
$ws.Cells.Item(4, 8).Value = 44881.449103900464
$ws.Cells.Item(4, 9).Value = 44881.450163761576

# Row 5
$ws.Cells.Item(5, 2).Value = "dfgdsg"
$ws.Cells.Item(5, 3).Value = "'+998908171355"
$ws.Cells.Item(5, 4).Value = "STORE_OWNER: bnjhbhjb COURIER: sdfsdfds"
$ws.Cells.Item(5, 5).Value = "REJECTED_DELIVERING"
$ws.Cells.Item(5, 7).Value = 1222
$ws.Cells.Item(5, 8).Value = 44881.44031289352
$ws.Cells.Item(5, 9).Value = 44881.4628338426

# Row 6
$ws.Cells.Item(6, 2).Value = "bfdfbs"
$ws.Cells.Item(6, 3).Value = "'+998908171355"
$ws.Cells.Item(6, 4).Value = "STORE_OWNER: nhbdfbhvbf COURIER: dgdfg"
$ws.Cells.Item(6, 5).Value = "PENDING"
$ws.Cells.Item(6, 7).Value = 12222
$ws.Cells.Item(6, 8).Value = 44881.44031289352
$ws.Cells.Item(6, 9).Value = 44881.46646645833

# Row 7
$ws.Cells.Item(7, 2).Value = "sdgdfsgq"
$ws.Cells.Item(7, 3).Value = "'+998908171355"
$ws.Cells.Item(7, 4).Value = "STORE_OWNER: fdsfsdf COURIER: sfsafs"
$ws.Cells.Item(7, 5).Value = "SOLD"
$ws.Cells.Item(7, 7).Value = 111111
$ws.Cells.Item(7, 8).Value = 44881.44031289352
$ws.Cells.Item(7, 9).Value = 44881.521820486116
